$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the I (I0) and J (IF) columns for rows 2-43
$iValues = @(5,4,6,1,3,4,7,6,6,5,6,5,10,8,7,9,4,7,8,8,4,11,10,6,9,7,3,8,9,8,9,6,6,7,7,5,7,9,4,8,9,5)
$jValues = @(5,5,7,3,4,5,8,7,6,5,8,6,10,8,7,9,5,8,8,9,6,11,11,7,9,7,4,9,9,8,9,7,7,7,8,7,8,9,5,8,9,5)

for ($r = 2; $r -le 43; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
